$wb = $excel.ActiveWorkbook

# Rename "Paineis DARQ" -> "PAINEIS DARQ"
$ws1 = $wb.Sheets.Item("Paineis DARQ")
$ws1.Name = "PAINEIS DARQ"

# Rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
$ws2 = $wb.Sheets.Item("Recolhimento x Eliminacao")
$ws2.Name = "RECOLHIMENTO X ELIMINAÇÃO"

# Delete the "Desarquivamentos Pendentes" sheet entirely
$excel.DisplayAlerts = $false
$ws3 = $wb.Sheets.Item("Desarquivamentos Pendentes")
$ws3.Delete()
$excel.DisplayAlerts = $true
